# Adicionando Bundesliga + BR
# Insert a new alias row for "Darmstadt" (sigla D98) right after the
# existing "Darmstadt 98" row (row 57), shifting everything below it down
# by one row, then update the view/selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 57 (pushes old row 57.. down to 58..)
$ws.Rows.Item(57).Insert()

# Populate the newly inserted row with the new alias
$ws.Cells.Item(57, 1).Value = "Darmstadt"
$ws.Cells.Item(57, 2).Value = "D98"

# Resize the AutoFilter to include the newly inserted row. Toggling it off
# then re-applying over the new range avoids relying on auto-expand.
$ws.AutoFilterMode = $false
$ws.Range("A1:B70").AutoFilter()

# Match the saved view position / selection from the edit
$ws.Application.ActiveWindow.ScrollRow = 53
$ws.Range("B57").Select()
